# Auto-generated Excel COM-interop script to apply cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.255.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.63%  "
$ws.Range("D3").Value = "'2.979.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.37%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'544.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").Value = "'130.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.19%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "'2.974.29"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.35%  "
$ws.Range("E9").Value = "  -2.56%  "
$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D10").Value = "'5.92"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.14%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.143"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.98%  "
$ws.Range("E12").Value = "  -3.96%  "
$ws.Range("D13").Value = "'0.0000216"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.61%  "
$ws.Range("D14").Value = "'33.42"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.79%  "
$ws.Range("D15").Value = "'3.460.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.48%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "'61.203.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.70%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.109"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.80%  "
$ws.Range("D18").Value = "'2.974.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.60%  "
$ws.Range("D19").Value = "'6.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.94%  "
$ws.Range("D20").Value = "'470.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.08%  "
$ws.Range("D21").Value = "'12.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.73%  "
$ws.Range("D22").Value = "'0.658"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.24%  "
$ws.Range("D23").Value = "'6.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.87%  "
$ws.Range("D24").Value = "'79.33"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("D25").Value = "'11.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.39%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -1.57%  "
$ws.Range("D28").Value = "'7.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.62%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").Value = "'1.87"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.25%  "
$ws.Range("D31").Value = "'25.30"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.57%  "
$ws.Range("D32").Value = "'1.12"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.66%  "
$ws.Range("E33").Value = "  -2.49%  "
$ws.Range("D34").Value = "'5.38"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.87%  "
$ws.Range("D35").Value = "'54.41"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.00%  "
$ws.Range("D36").Value = "'5.80"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.64%  "
$ws.Range("D37").Value = "'443.04"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -10.27%  "
$ws.Range("D38").Value = "'3.114.50"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.10%  "
$ws.Range("D39").Value = "'0.0782"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.34%  "
$ws.Range("D40").Value = "'0.0374"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.61%  "
$ws.Range("E41").Value = "  -2.89%  "
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").Value = "'2.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -13.24%  "
$ws.Range("D45").Value = "'25.06"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.73%  "
$ws.Range("D46").Value = "'0.238"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.74%  "
$ws.Range("E47").Value = "  -2.60%  "
$ws.Range("E49").Value = "  -6.61%  "
$ws.Range("D50").Value = "'113.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -9.33%  "
$ws.Range("D51").Value = "'0.0₃0475"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -10.24%  "
